$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 6-10 (column A) hold the restriction formulas; column B gets a new
# label naming which product each restriction refers to. Clear any stale
# formatting first so the new label cells pick up a clean "italic only"
# style instead of inheriting the old centered style that used to live on
# B8:B10 (when they were blank placeholder cells).
$labels = $ws.Range("B6:B10")
$labels.ClearFormats()

# Assign values in this order so the shared-string table is built up the
# same way the original author's workbook has it (diesel, gasolina,
# lubrificante, combustível para jatos, then óleo cru last).
$ws.Range("B7").Value = "diesel"
$ws.Range("B8").Value = "gasolina"
$ws.Range("B9").Value = "lubrificante"
$ws.Range("B10").Value = "combustível para jatos"
$ws.Range("B6").Value = "óleo cru"

# Style the new labels in italics, matching the rest of the sheet's label
# formatting conventions.
$labels.Font.Italic = $true

# Move the active selection to B11, where the user's cursor ended up.
$ws.Range("B11").Select()
